# mise en commentaires des paramètres de recherche qu'on n'utilise plus
#
# Refresh the "Date" metadata value and add a new "Jurisdiction" row
# (currently left blank) to the Metadata sheet, right after "Contact"
# and before "Description" - this pushes Description/Purpose/Copyright/
# Immutable down by one row.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Refresh the Date value (row 8, column B)
$ws1.Range("B8").Value = "2024-07-01T07:50:29+00:00"

# Insert a new row before the old "Description" row (row 11)
$ws1.Rows.Item(11).Insert()
$ws1.Range("A11").Value = "Jurisdiction"
$ws1.Range("B11").Value = ""

# Give the new row the same look (font/fill/borders) as its neighbours
$ws1.Range("A12:B12").Copy()
$ws1.Range("A11:B11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# PasteSpecial only touches formatting, but make sure the values are
# still exactly what we want
$ws1.Range("A11").Value = "Jurisdiction"
$ws1.Range("B11").Value = ""
